# Apply fixes to the 3.3V and 12V schematic calculator sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "3.3V" ---
$ws33 = $wb.Worksheets.Item("3.3V")
$ws33.Activate()
$ws33.Range("B29").Value = 600000
$ws33.Range("B29").Select()

# --- Sheet "12V" ---
$ws12 = $wb.Worksheets.Item("12V")
$ws12.Activate()
$ws12.Range("B17").Formula = "=10000+470"
$ws12.Range("B29").Value = 500000
$ws12.Range("B29").Select()

$wb.Save()
